# Fix error handling for bad file extensions
# Update sample snippet data rows with corrected sampling values
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 49
$ws.Range("B2").Value = "5:50 AM"
$ws.Range("C2").Value = 1467.18
$ws.Range("D2").Value = 1497.18
$ws.Range("E2").Value = 24.93
$ws.Range("A3").Value = 91
$ws.Range("B3").Value = "6:11 AM"
$ws.Range("C3").Value = 2725.88
$ws.Range("D3").Value = 2755.88
$ws.Range("E3").Value = 22.16
$ws.Range("A4").Value = 100
$ws.Range("B4").Value = "6:15 AM"
$ws.Range("C4").Value = 2986.3175
$ws.Range("D4").Value = 3016.3175
$ws.Range("E4").Value = 30.9975
$ws.Range("A5").Value = 123
$ws.Range("B5").Value = "6:27 AM"
$ws.Range("C5").Value = 3687.43
$ws.Range("D5").Value = 3717.43
$ws.Range("E5").Value = 25.2
$ws.Range("A6").Value = 151
$ws.Range("B6").Value = "6:40 AM"
$ws.Range("C6").Value = 4506.15
$ws.Range("D6").Value = 4536.15
$ws.Range("E6").Value = 33.81
$ws.Range("A7").Value = 174
$ws.Range("B7").Value = "6:52 AM"
$ws.Range("C7").Value = 5193.04
$ws.Range("D7").Value = 5223.04
$ws.Range("E7").Value = 71.97
$ws.Range("A11").Value = 251
$ws.Range("B11").Value = "7:30 AM"
$ws.Range("C11").Value = 7516.66
$ws.Range("D11").Value = 7546.66
$ws.Range("E11").Value = 27.63
$ws.Range("A15").Value = 993
$ws.Range("B15").Value = "1:41 PM"
$ws.Range("C15").Value = 29775.52
$ws.Range("D15").Value = 29805.52
$ws.Range("E15").Value = 41.35
$ws.Range("A20").Value = 1141
$ws.Range("B20").Value = "2:55 PM"
$ws.Range("C20").Value = 34219.87
$ws.Range("D20").Value = 34249.87
$ws.Range("E20").Value = 21.48
$ws.Range("A21").Value = 1180
$ws.Range("B21").Value = "3:15 PM"
$ws.Range("C21").Value = 35378.663333
$ws.Range("D21").Value = 35408.663333
$ws.Range("E21").Value = 42.293333
$ws.Range("A22").Value = 1190
$ws.Range("B22").Value = "3:20 PM"
$ws.Range("C22").Value = 35678.675
$ws.Range("D22").Value = 35708.675
$ws.Range("E22").Value = 40.53
$ws.Range("A23").Value = 1196
$ws.Range("B23").Value = "3:23 PM"
$ws.Range("C23").Value = 35858.715
$ws.Range("D23").Value = 35888.715
$ws.Range("E23").Value = 56.205
$ws.Range("A24").Value = 1202
$ws.Range("B24").Value = "3:26 PM"
$ws.Range("C24").Value = 36041.01
$ws.Range("D24").Value = 36071.01
$ws.Range("E24").Value = 24.126667
$ws.Range("A32").Value = 1299
$ws.Range("B32").Value = "4:14 PM"
$ws.Range("C32").Value = 38945.82
$ws.Range("D32").Value = 38975.82
$ws.Range("E32").Value = 32.71
$ws.Range("A33").Value = 1308
$ws.Range("B33").Value = "4:19 PM"
$ws.Range("C33").Value = 39217.913333
$ws.Range("D33").Value = 39247.913333
$ws.Range("E33").Value = 29.03
$ws.Range("A34").Value = 1315
$ws.Range("B34").Value = "4:22 PM"
$ws.Range("C34").Value = 39427.435
$ws.Range("D34").Value = 39457.435
$ws.Range("E34").Value = 98.93000000000001
$ws.Range("A35").Value = 1323
$ws.Range("B35").Value = "4:26 PM"
$ws.Range("C35").Value = 39677.255
$ws.Range("D35").Value = 39707.255
$ws.Range("E35").Value = 57.735
$ws.Range("A36").Value = 1335
$ws.Range("B36").Value = "4:32 PM"
$ws.Range("C36").Value = 40027.72
$ws.Range("D36").Value = 40057.72
$ws.Range("E36").Value = 47.705
$ws.Range("A37").Value = 1343
$ws.Range("B37").Value = "4:36 PM"
$ws.Range("C37").Value = 40265.54
$ws.Range("D37").Value = 40295.54
$ws.Range("E37").Value = 30.27
$ws.Range("A38").Value = 1354
$ws.Range("B38").Value = "4:42 PM"
$ws.Range("C38").Value = 40599.6
$ws.Range("D38").Value = 40629.6
$ws.Range("E38").Value = 20.99
$ws.Range("A39").Value = 1362
$ws.Range("B39").Value = "4:46 PM"
$ws.Range("C39").Value = 40838.935
$ws.Range("D39").Value = 40868.935
$ws.Range("E39").Value = 104.33
$ws.Range("A40").Value = 1372
$ws.Range("B40").Value = "4:51 PM"
$ws.Range("C40").Value = 41141.465
$ws.Range("D40").Value = 41171.465
$ws.Range("E40").Value = 98.08
$ws.Range("A41").Value = 1390
$ws.Range("B41").Value = "5:00 PM"
$ws.Range("C41").Value = 41680.48
$ws.Range("D41").Value = 41710.48
$ws.Range("E41").Value = 38.76
$ws.Range("A42").Value = 1408
$ws.Range("B42").Value = "5:09 PM"
$ws.Range("C42").Value = 42222.19
$ws.Range("D42").Value = 42252.19
$ws.Range("E42").Value = 25.425
$ws.Range("A43").Value = 1418
$ws.Range("B43").Value = "5:14 PM"
$ws.Range("C43").Value = 42521.595
$ws.Range("D43").Value = 42551.595
$ws.Range("E43").Value = 25.965
$ws.Range("A44").Value = 1427
$ws.Range("B44").Value = "5:19 PM"
$ws.Range("C44").Value = 42809.97
$ws.Range("D44").Value = 42839.97
$ws.Range("E44").Value = 161.82
$ws.Range("A45").Value = 1441
$ws.Range("B45").Value = "5:26 PM"
$ws.Range("C45").Value = 43229.93
$ws.Range("D45").Value = 43259.93
$ws.Range("E45").Value = 24.11
$ws.Range("A51").Value = 1551
$ws.Range("B51").Value = "6:20 PM"
$ws.Range("C51").Value = 46502.28
$ws.Range("D51").Value = 46532.28
$ws.Range("E51").Value = 39.27
